# Updates cryptos list values per the Fri Sep  1 03:10:57 UTC 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D sometimes holds purely numeric-looking price strings (e.g. "217.77").
# Excel auto-converts those to real numbers on assignment, which would flip the
# cell from a text cell to a numeric one. Mark them as Text first so the original
# "price as string" representation is preserved, then restore the default style so
# no stray per-cell formatting is left behind.
$numericLookingPriceCells = @(
    "D5", "D6", "D8", "D9", "D10", "D11", "D13", "D15", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D30", "D32", "D33", "D34", "D36", "D37", "D40", "D41", "D42", "D43", "D47", "D48", "D49", "D50", "D51"
)
foreach ($cellRef in $numericLookingPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '26.186.11'
$ws.Range('D3').Value = '1.657.45'
$ws.Range('E3').Value = '  -2.82%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').Value = '217.77'
$ws.Range('E5').Value = '  -2.72%  '
$ws.Range('D6').Value = '0.5143'
$ws.Range('E6').Value = '  -3.10%  '
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('D8').Value = '0.2582'
$ws.Range('E8').Value = '  -2.94%  '
$ws.Range('D9').Value = '0.06433'
$ws.Range('E9').Value = '  -2.27%  '
$ws.Range('D10').Value = '19.93'
$ws.Range('E10').Value = '  -3.83%  '
$ws.Range('D11').Value = '0.07813'
$ws.Range('E11').Value = '  +2.37%  '
$ws.Range('D12').Value = '1.662.82'
$ws.Range('E12').Value = '  -2.62%  '
$ws.Range('D13').Value = '4.290'
$ws.Range('E13').Value = '  -4.59%  '
$ws.Range('D14').Value = '1.885.66'
$ws.Range('E14').Value = '  -2.83%  '
$ws.Range('D15').Value = '0.5542'
$ws.Range('E15').Value = '  -4.08%  '
$ws.Range('D16').Value = '0.0₅8061'
$ws.Range('E16').Value = '  -1.02%  '
$ws.Range('D17').Value = '64.20'
$ws.Range('E17').Value = '  -5.03%  '
$ws.Range('D18').Value = '26.218.79'
$ws.Range('E18').Value = '  -3.99%  '
$ws.Range('D19').Value = '211.82'
$ws.Range('E19').Value = '  -1.60%  '
$ws.Range('D20').Value = '1.005'
$ws.Range('E20').Value = '  +0.10%  '
$ws.Range('D21').Value = '4.423'
$ws.Range('E21').Value = '  -4.05%  '
$ws.Range('D22').Value = '10.02'
$ws.Range('E22').Value = '  -3.34%  '
$ws.Range('D23').Value = '5.962'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').Value = '1.006'
$ws.Range('E24').Value = '  +0.13%  '
$ws.Range('D25').Value = '143.34'
$ws.Range('E25').Value = '  -0.65%  '
$ws.Range('D26').Value = '1.755'
$ws.Range('E26').Value = '  +3.37%  '
$ws.Range('D27').Value = '0.1163'
$ws.Range('E27').Value = '  -3.16%  '
$ws.Range('D28').Value = '6.968'
$ws.Range('E28').Value = '  -3.32%  '
$ws.Range('E29').Value = '  -1.88%  '
$ws.Range('D30').Value = '0.05230'
$ws.Range('E30').Value = '  -2.61%  '
$ws.Range('E31').Value = '  -2.55%  '
$ws.Range('D32').Value = '3.365'
$ws.Range('E32').Value = '  -3.01%  '
$ws.Range('D33').Value = '3.217'
$ws.Range('E33').Value = '  -5.45%  '
$ws.Range('D34').Value = '1.568'
$ws.Range('E34').Value = '  -4.74%  '
$ws.Range('E35').Value = '  -3.75%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = '2.371'
$ws.Range('E36').Value = '  -1.82%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').Value = '0.9292'
$ws.Range('E37').Value = '  -1.77%  '
$ws.Range('D38').Value = '1.173.23'
$ws.Range('E38').Value = '  +12.65%  '
$ws.Range('E39').Value = '  -1.99%  '
$ws.Range('D40').Value = '0.01592'
$ws.Range('E40').Value = '  -2.15%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').Value = '1.005'
$ws.Range('E41').Value = '  +0.12%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = '0.8443'
$ws.Range('E42').Value = '  +0.42%  '
$ws.Range('D43').Value = '5.673'
$ws.Range('E43').Value = '  -1.83%  '
$ws.Range('E44').Value = '  -0.67%  '
$ws.Range('D45').Value = '1.796.23'
$ws.Range('E45').Value = '  -2.84%  '
$ws.Range('D46').Value = '0.0₈115'
$ws.Range('E46').Value = '  +2.79%  '
$ws.Range('D47').Value = '0.4539'
$ws.Range('E47').Value = '  +0.42%  '
$ws.Range('D48').Value = '55.92'
$ws.Range('E48').Value = '  -3.13%  '
$ws.Range('D49').Value = '1.002'
$ws.Range('E49').Value = '  -0.21%  '
$ws.Range('D50').Value = '7.857'
$ws.Range('E50').Value = '  -2.24%  '
$ws.Range('D51').Value = '0.05057'
$ws.Range('E51').Value = '  -3.26%  '

foreach ($cellRef in $numericLookingPriceCells) {
    $ws.Range($cellRef).Style = "Normal"
}
